$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "The applicant has been given a score of 81.79. The applicant’s projects demonstrate significant expertise in NLP, computer vision, and machine learning, which aligns with the requirements of the job description. The 'LLMGuard' project exhibits the ability to work with language models like GPT-2 and BERT, which is directly relevant to the job role. Additionally, the 'CloudPhysician's Vital Extraction Challenge' project showcases advanced skills in computer vision and image processing. The applicant’s proficiency in Python, PyTorch, and ReactJS is evident from the projects, further indicating strong alignment with the required skills for the job. However, the absence of explicit NLP project experience and the use of JavaScript, not mentioned in the required skills, may have influenced the score given to the applicant."

$ws.Range("E2").Value = "The candidate has expressed a deep interest in AI and a willingness to adapt to new cultures, which aligns with the soft skills required by the company. The applicant's response also indicates a keenness to contribute to cutting-edge technologies and a desire to work in Japan for professional development. Moreover, the emphasis on teamwork, adaptability, and problem-solving in overcoming challenges further reflects positive attributes. However, the preference to not work alone and the need to improve individual working skills slightly impact the overall score, leading to a rating of 4 out of 5."

$ws.Range("C3").Value = "The applicant has been involved in projects that align very closely with the job description. The 'Cloudphysician's Vital Extraction Challenge' and 'LLMGuard' projects demonstrate strong skills in Computer Vision, Image Processing, NLP, and Deep Learning, using frameworks like PyTorch and Python - all of which are essential for the NLP Engineer role. The 'Multi Model Data Analysis for Annotation of Human Activities' project also showcases proficiency in Python, PyTorch, and Tensorflow, indicating a good fit for the position. Additionally, the applicant's skills in Django Rest Framework and ReactJS from the 'LLMGuard' and 'Website for the Literature Society of the college' projects display versatility and experience in web development, which can be valuable for prototyping in NLP and computer vision. However, the applicant does not explicitly mention NLP experience in the projects, which might have slightly affected the score."

$ws.Range("E3").Value = "The applicant has demonstrated strong communication skills, adaptability to new cultures, and a willingness to learn and contribute to cutting-edge AI/ML technologies. The applicant's positive attitude towards working in Japan, openness to learning a new language, and expressed interest in being part of Japan's work culture reflects a high level of adaptability and enthusiasm. Additionally, the applicant's emphasis on teamwork, open communication, problem-solving, and the ability to bridge technical and creative aspects are aligned with the company's soft skills requirements. Overall, the applicant's eagerness to contribute, learn, and adapt indicates a strong fit for the company, warranting a score of 4."
